# Auto-generated update of cryptos price/volume columns
# Applies latest scraped price (D) and 1h volume change (E) values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.615.17'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.643.84'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = '''216.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("E6").Value = '  +0.55%  '
$ws.Range("D7").Value = '''1.01'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").Value = '''0.0626'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").Value = '''19.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("D11").Value = '''0.0844'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '1.873.98'
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("E13").Value = '  +3.15%  '
$ws.Range("D14").Value = '1.612.50'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").Value = '''0.533'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("D16").Value = '''65.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.16%  '
$ws.Range("D17").Value = '26.673.06'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = '0.0₃0750'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").Value = '''218.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").Value = '''9.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("E24").Value = '  +9.74%  '
$ws.Range("D25").Value = '''146.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.29%  '
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").Value = '''7.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.01%  '
$ws.Range("D29").Value = '''15.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.46%  '
$ws.Range("D30").Value = '''0.0517'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("E31").Value = '  +1.29%  '
$ws.Range("D32").Value = '''3.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("E33").Value = '  +2.74%  '
$ws.Range("D34").Value = '1.274.11'
$ws.Range("E34").Value = '  +5.19%  '
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("D36").Value = '''0.0183'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.10%  '
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = '''0.526'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.48%  '
$ws.Range("E39").Value = '  +1.97%  '
$ws.Range("D40").Value = '''1.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").Value = '''0.806'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("E42").Value = '  -1.94%  '
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("D44").Value = '1.785.47'
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("D45").Value = '''93.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").Value = '''59.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.93%  '
$ws.Range("D47").Value = '''1.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.48%  '
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = '''7.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("D50").Value = '''0.0976'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.59%  '
$ws.Range("E51").Value = '  -0.53%  '
